$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-03-20 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-21 Friday", 2) | Out-Null
$d.Content.Find.Execute("544×8=4352", $true, $false, $false, $false, $false, $true, 1, $false, "996×9=8964", 2) | Out-Null
$d.Content.Find.Execute("863×6=5178", $true, $false, $false, $false, $false, $true, 1, $false, "150×7=1050", 2) | Out-Null
$d.Content.Find.Execute("438×5=2190", $true, $false, $false, $false, $false, $true, 1, $false, "436×8=3488", 2) | Out-Null
$d.Content.Find.Execute("768×6=4608", $true, $false, $false, $false, $false, $true, 1, $false, "209×5=1045", 2) | Out-Null
$d.Content.Find.Execute("330×8=2640", $true, $false, $false, $false, $false, $true, 1, $false, "390×8=3120", 2) | Out-Null
$d.Content.Find.Execute("512×2=1024", $true, $false, $false, $false, $false, $true, 1, $false, "237×6=1422", 2) | Out-Null
$d.Content.Find.Execute("754×9=6786", $true, $false, $false, $false, $false, $true, 1, $false, "637×8=5096", 2) | Out-Null
$d.Content.Find.Execute("870×6=5220", $true, $false, $false, $false, $false, $true, 1, $false, "303×4=1212", 2) | Out-Null
$d.Content.Find.Execute("205×4=820", $true, $false, $false, $false, $false, $true, 1, $false, "229×7=1603", 2) | Out-Null
$d.Content.Find.Execute("281×3=843", $true, $false, $false, $false, $false, $true, 1, $false, "863×7=6041", 2) | Out-Null
$d.Content.Find.Execute("825×8=6600", $true, $false, $false, $false, $false, $true, 1, $false, "848×6=5088", 2) | Out-Null
$d.Content.Find.Execute("224×2=448", $true, $false, $false, $false, $false, $true, 1, $false, "306×2=612", 2) | Out-Null
$d.Content.Find.Execute("595×5=2975", $true, $false, $false, $false, $false, $true, 1, $false, "630×9=5670", 2) | Out-Null
$d.Content.Find.Execute("302×6=1812", $true, $false, $false, $false, $false, $true, 1, $false, "282×3=846", 2) | Out-Null
$d.Content.Find.Execute("982×6=5892", $true, $false, $false, $false, $false, $true, 1, $false, "305×7=2135", 2) | Out-Null
$d.Content.Find.Execute("332×9=2988", $true, $false, $false, $false, $false, $true, 1, $false, "966×8=7728", 2) | Out-Null
$d.Content.Find.Execute("257×5=1285", $true, $false, $false, $false, $false, $true, 1, $false, "561×5=2805", 2) | Out-Null
$d.Content.Find.Execute("352×7=2464", $true, $false, $false, $false, $false, $true, 1, $false, "550×5=2750", 2) | Out-Null
$d.Content.Find.Execute("514×3=1542", $true, $false, $false, $false, $false, $true, 1, $false, "779×8=6232", 2) | Out-Null
$d.Content.Find.Execute("878×2=1756", $true, $false, $false, $false, $false, $true, 1, $false, "159×4=636", 2) | Out-Null
$d.Content.Find.Execute("905×6=5430", $true, $false, $false, $false, $false, $true, 1, $false, "275×6=1650", 2) | Out-Null
$d.Content.Find.Execute("591×9=5319", $true, $false, $false, $false, $false, $true, 1, $false, "880×8=7040", 2) | Out-Null
$d.Content.Find.Execute("932×5=4660", $true, $false, $false, $false, $false, $true, 1, $false, "183×3=549", 2) | Out-Null
$d.Content.Find.Execute("285×7=1995", $true, $false, $false, $false, $false, $true, 1, $false, "879×4=3516", 2) | Out-Null
$d.Content.Find.Execute("340×4=1360", $true, $false, $false, $false, $false, $true, 1, $false, "737×4=2948", 2) | Out-Null
